$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the test case in row 2 from "Test Case Happy Path" to "Test Case 1"
$ws.Range("A2").Value = "Test Case 1"

# Move/record the active selection to A8 (was D13)
[void]$ws.Range("A8").Select()
